$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Reln"
$ws.Cells.Item(2, 3).Value = "Itga3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 0.09865600000000001
$ws.Cells.Item(2, 8).Value = 0.295968
$ws.Cells.Item(2, 9).Value = 0.01009304870291239
$ws.Cells.Item(2, 10).Value = 0.01488643315542961
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 12).Value = 1.0
$ws.Cells.Item(2, 13).Value = 11.451657
$ws.Cells.Item(2, 14).Value = 34.354971
$ws.Cells.Item(2, 15).Value = 0.6845114669615147
$ws.Cells.Item(2, 16).Value = 0.752353342656931
$ws.Cells.Item(2, 17).Value = 1.129774672992
$ws.Cells.Item(2, 18).Value = 10.167972056928
$ws.Cells.Item(2, 19).Value = 0.006908807573744576
$ws.Cells.Item(2, 20).Value = 0.01119985774472643

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Reln"
$ws.Cells.Item(3, 3).Value = "Itga3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 0.09865600000000001
$ws.Cells.Item(3, 8).Value = 0.295968
$ws.Cells.Item(3, 9).Value = 0.01009304870291239
$ws.Cells.Item(3, 10).Value = 0.01488643315542961
$ws.Cells.Item(3, 11).Value = 3.0
$ws.Cells.Item(3, 12).Value = 1.0
$ws.Cells.Item(3, 13).Value = 0.5240143333333332
$ws.Cells.Item(3, 14).Value = 1.572043
$ws.Cells.Item(3, 15).Value = 0.03132243831777883
$ws.Cells.Item(3, 16).Value = 0.03442680262633404
$ws.Cells.Item(3, 17).Value = 0.05169715806933333
$ws.Cells.Item(3, 18).Value = 0.465274422624
$ws.Cells.Item(3, 19).Value = 0.000316138895435311
$ws.Cells.Item(3, 20).Value = 0.0005124922960520901

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Reln"
$ws.Cells.Item(4, 3).Value = "Itga3"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 0.09865600000000001
$ws.Cells.Item(4, 8).Value = 0.295968
$ws.Cells.Item(4, 9).Value = 0.01009304870291239
$ws.Cells.Item(4, 10).Value = 0.01488643315542961
$ws.Cells.Item(4, 11).Value = 1.0
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.123396
$ws.Cells.Item(4, 14).Value = 0.370188
$ws.Cells.Item(4, 15).Value = 0.007375873812600488
$ws.Cells.Item(4, 16).Value = 0.00810689606495328
$ws.Cells.Item(4, 17).Value = 0.012173755776
$ws.Cells.Item(4, 18).Value = 0.109563801984
$ws.Cells.Item(4, 19).Value = 0.00007444505361711285
$ws.Cells.Item(4, 20).Value = 0.0001206827663689423

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Reln"
$ws.Cells.Item(5, 3).Value = "Itga3"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 0.09865600000000001
$ws.Cells.Item(5, 8).Value = 0.295968
$ws.Cells.Item(5, 9).Value = 0.01009304870291239
$ws.Cells.Item(5, 10).Value = 0.01488643315542961
$ws.Cells.Item(5, 11).Value = 2.0
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.104921
$ws.Cells.Item(5, 14).Value = 0.314763
$ws.Cells.Item(5, 15).Value = 0.006271548966675223
$ws.Cells.Item(5, 16).Value = 0.006893121673562862
$ws.Cells.Item(5, 17).Value = 0.010351086176
$ws.Cells.Item(5, 18).Value = 0.093159775584
$ws.Cells.Item(5, 19).Value = 0.00006329904916335292
$ws.Cells.Item(5, 20).Value = 0.0001026139950257366

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Reln"
$ws.Cells.Item(6, 3).Value = "Itga3"
$ws.Cells.Item(6, 4).Value = "sCs"
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 0.09865600000000001
$ws.Cells.Item(6, 8).Value = 0.295968
$ws.Cells.Item(6, 9).Value = 0.01009304870291239
$ws.Cells.Item(6, 10).Value = 0.01488643315542961
$ws.Cells.Item(6, 11).Value = 2.0
$ws.Cells.Item(6, 12).Value = 1.0
$ws.Cells.Item(6, 13).Value = 4.525690500000001
$ws.Cells.Item(6, 14).Value = 9.051381000000001
$ws.Cells.Item(6, 15).Value = 0.2705186719414309
$ws.Cells.Item(6, 16).Value = 0.1982198369782188
$ws.Cells.Item(6, 17).Value = 0.4464865219680001
$ws.Cells.Item(6, 18).Value = 2.678919131808001
$ws.Cells.Item(6, 19).Value = 0.002730358130952043
$ws.Cells.Item(6, 20).Value = 0.002950786353256408

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Reln"
$ws.Cells.Item(7, 3).Value = "Itga3"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 0.2337746666666667
$ws.Cells.Item(7, 8).Value = 0.7013240000000001
$ws.Cells.Item(7, 9).Value = 0.02391642775070728
$ws.Cells.Item(7, 10).Value = 0.03527480283780177
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 12).Value = 1.0
$ws.Cells.Item(7, 13).Value = 11.451657
$ws.Cells.Item(7, 14).Value = 34.354971
$ws.Cells.Item(7, 15).Value = 0.6845114669615147
$ws.Cells.Item(7, 16).Value = 0.752353342656931
$ws.Cells.Item(7, 17).Value = 2.677107297956
$ws.Cells.Item(7, 18).Value = 24.093965681604
$ws.Cells.Item(7, 19).Value = 0.01637106904411572
$ws.Cells.Item(7, 20).Value = 0.02653911582658436

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Reln"
$ws.Cells.Item(8, 3).Value = "Itga3"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 0.2337746666666667
$ws.Cells.Item(8, 8).Value = 0.7013240000000001
$ws.Cells.Item(8, 9).Value = 0.02391642775070728
$ws.Cells.Item(8, 10).Value = 0.03527480283780177
$ws.Cells.Item(8, 11).Value = 3.0
$ws.Cells.Item(8, 12).Value = 1.0
$ws.Cells.Item(8, 13).Value = 0.5240143333333332
$ws.Cells.Item(8, 14).Value = 1.572043
$ws.Cells.Item(8, 15).Value = 0.03132243831777883
$ws.Cells.Item(8, 16).Value = 0.03442680262633404
$ws.Cells.Item(8, 17).Value = 0.1225012761035555
$ws.Cells.Item(8, 18).Value = 1.102511484932
$ws.Cells.Item(8, 19).Value = 0.0007491208330031425
$ws.Cells.Item(8, 20).Value = 0.001214398674979849

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Reln"
$ws.Cells.Item(9, 3).Value = "Itga3"
$ws.Cells.Item(9, 4).Value = "M1"
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 0.2337746666666667
$ws.Cells.Item(9, 8).Value = 0.7013240000000001
$ws.Cells.Item(9, 9).Value = 0.02391642775070728
$ws.Cells.Item(9, 10).Value = 0.03527480283780177
$ws.Cells.Item(9, 11).Value = 1.0
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.123396
$ws.Cells.Item(9, 14).Value = 0.370188
$ws.Cells.Item(9, 15).Value = 0.007375873812600488
$ws.Cells.Item(9, 16).Value = 0.00810689606495328
$ws.Cells.Item(9, 17).Value = 0.028846858768
$ws.Cells.Item(9, 18).Value = 0.259621728912
$ws.Cells.Item(9, 19).Value = 0.0001764045531373934
$ws.Cells.Item(9, 20).Value = 0.0002859691603177779

# Row 10
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Reln"
$ws.Cells.Item(10, 3).Value = "Itga3"
$ws.Cells.Item(10, 4).Value = "M2"
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 0.2337746666666667
$ws.Cells.Item(10, 8).Value = 0.7013240000000001
$ws.Cells.Item(10, 9).Value = 0.02391642775070728
$ws.Cells.Item(10, 10).Value = 0.03527480283780177
$ws.Cells.Item(10, 11).Value = 2.0
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.104921
$ws.Cells.Item(10, 14).Value = 0.314763
$ws.Cells.Item(10, 15).Value = 0.006271548966675223
$ws.Cells.Item(10, 16).Value = 0.006893121673562862
$ws.Cells.Item(10, 17).Value = 0.02452787180133334
$ws.Cells.Item(10, 18).Value = 0.220750846212
$ws.Cells.Item(10, 19).Value = 0.0001499930477465109
$ws.Cells.Item(10, 20).Value = 0.0002431535079719081

# Row 11
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Reln"
$ws.Cells.Item(11, 3).Value = "Itga3"
$ws.Cells.Item(11, 4).Value = "sCs"
$ws.Cells.Item(11, 5).Value = 3.0
$ws.Cells.Item(11, 6).Value = 1.0
$ws.Cells.Item(11, 7).Value = 0.2337746666666667
$ws.Cells.Item(11, 8).Value = 0.7013240000000001
$ws.Cells.Item(11, 9).Value = 0.02391642775070728
$ws.Cells.Item(11, 10).Value = 0.03527480283780177
$ws.Cells.Item(11, 11).Value = 2.0
$ws.Cells.Item(11, 12).Value = 1.0
$ws.Cells.Item(11, 13).Value = 4.525690500000001
$ws.Cells.Item(11, 14).Value = 9.051381000000001
$ws.Cells.Item(11, 15).Value = 0.2705186719414309
$ws.Cells.Item(11, 16).Value = 0.1982198369782188
$ws.Cells.Item(11, 17).Value = 1.057991788074
$ws.Cells.Item(11, 18).Value = 6.347950728444001
$ws.Cells.Item(11, 19).Value = 0.006469840272704517
$ws.Cells.Item(11, 20).Value = 0.006992165667947877

# Row 12
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Reln"
$ws.Cells.Item(12, 3).Value = "Itga3"
$ws.Cells.Item(12, 4).Value = "ECs"
$ws.Cells.Item(12, 5).Value = 2.0
$ws.Cells.Item(12, 6).Value = 1.0
$ws.Cells.Item(12, 7).Value = 9.4422175
$ws.Cells.Item(12, 8).Value = 18.884435
$ws.Cells.Item(12, 9).Value = 0.9659905235463803
$ws.Cells.Item(12, 10).Value = 0.9498387640067686
$ws.Cells.Item(12, 11).Value = 3.0
$ws.Cells.Item(12, 12).Value = 1.0
$ws.Cells.Item(12, 13).Value = 11.451657
$ws.Cells.Item(12, 14).Value = 34.354971
$ws.Cells.Item(12, 15).Value = 0.6845114669615147
$ws.Cells.Item(12, 16).Value = 0.752353342656931
$ws.Cells.Item(12, 17).Value = 108.1290361293975
$ws.Cells.Item(12, 18).Value = 648.774216776385
$ws.Cells.Item(12, 19).Value = 0.6612315903436543
$ws.Cells.Item(12, 20).Value = 0.7146143690856201

# Row 13
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Reln"
$ws.Cells.Item(13, 3).Value = "Itga3"
$ws.Cells.Item(13, 4).Value = "FAPs"
$ws.Cells.Item(13, 5).Value = 2.0
$ws.Cells.Item(13, 6).Value = 1.0
$ws.Cells.Item(13, 7).Value = 9.4422175
$ws.Cells.Item(13, 8).Value = 18.884435
$ws.Cells.Item(13, 9).Value = 0.9659905235463803
$ws.Cells.Item(13, 10).Value = 0.9498387640067686
$ws.Cells.Item(13, 11).Value = 3.0
$ws.Cells.Item(13, 12).Value = 1.0
$ws.Cells.Item(13, 13).Value = 0.5240143333333332
$ws.Cells.Item(13, 14).Value = 1.572043
$ws.Cells.Item(13, 15).Value = 0.03132243831777883
$ws.Cells.Item(13, 16).Value = 0.03442680262633404
$ws.Cells.Item(13, 17).Value = 4.947857308450833
$ws.Cells.Item(13, 18).Value = 29.687143850705
$ws.Cells.Item(13, 19).Value = 0.03025717858934037
$ws.Cells.Item(13, 20).Value = 0.0326999116553021

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Reln"
$ws.Cells.Item(14, 3).Value = "Itga3"
$ws.Cells.Item(14, 4).Value = "M1"
$ws.Cells.Item(14, 5).Value = 2.0
$ws.Cells.Item(14, 6).Value = 1.0
$ws.Cells.Item(14, 7).Value = 9.4422175
$ws.Cells.Item(14, 8).Value = 18.884435
$ws.Cells.Item(14, 9).Value = 0.9659905235463803
$ws.Cells.Item(14, 10).Value = 0.9498387640067686
$ws.Cells.Item(14, 11).Value = 1.0
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.123396
$ws.Cells.Item(14, 14).Value = 0.370188
$ws.Cells.Item(14, 15).Value = 0.007375873812600488
$ws.Cells.Item(14, 16).Value = 0.00810689606495328
$ws.Cells.Item(14, 17).Value = 1.16513187063
$ws.Cells.Item(14, 18).Value = 6.990791223780001
$ws.Cells.Item(14, 19).Value = 0.007125024205845982
$ws.Cells.Item(14, 20).Value = 0.007700244138266559

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Reln"
$ws.Cells.Item(15, 3).Value = "Itga3"
$ws.Cells.Item(15, 4).Value = "M2"
$ws.Cells.Item(15, 5).Value = 2.0
$ws.Cells.Item(15, 6).Value = 1.0
$ws.Cells.Item(15, 7).Value = 9.4422175
$ws.Cells.Item(15, 8).Value = 18.884435
$ws.Cells.Item(15, 9).Value = 0.9659905235463803
$ws.Cells.Item(15, 10).Value = 0.9498387640067686
$ws.Cells.Item(15, 11).Value = 2.0
$ws.Cells.Item(15, 12).Value = 0.6666666666666666
$ws.Cells.Item(15, 13).Value = 0.104921
$ws.Cells.Item(15, 14).Value = 0.314763
$ws.Cells.Item(15, 15).Value = 0.006271548966675223
$ws.Cells.Item(15, 16).Value = 0.006893121673562862
$ws.Cells.Item(15, 17).Value = 0.9906869023175
$ws.Cells.Item(15, 18).Value = 5.944121413905
$ws.Cells.Item(15, 19).Value = 0.006058256869765358
$ws.Cells.Item(15, 20).Value = 0.006547354170565217

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Reln"
$ws.Cells.Item(16, 3).Value = "Itga3"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 2.0
$ws.Cells.Item(16, 6).Value = 1.0
$ws.Cells.Item(16, 7).Value = 9.4422175
$ws.Cells.Item(16, 8).Value = 18.884435
$ws.Cells.Item(16, 9).Value = 0.9659905235463803
$ws.Cells.Item(16, 10).Value = 0.9498387640067686
$ws.Cells.Item(16, 11).Value = 2.0
$ws.Cells.Item(16, 12).Value = 1.0
$ws.Cells.Item(16, 13).Value = 4.525690500000001
$ws.Cells.Item(16, 14).Value = 9.051381000000001
$ws.Cells.Item(16, 15).Value = 0.2705186719414309
$ws.Cells.Item(16, 16).Value = 0.1982198369782188
$ws.Cells.Item(16, 17).Value = 42.73255403868375
$ws.Cells.Item(16, 18).Value = 170.930216154735
$ws.Cells.Item(16, 19).Value = 0.2613184735377743
$ws.Cells.Item(16, 20).Value = 0.1882768849570146
